$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6491417538996131
$ws.Range("C2").Value = 0.256339892275065
$ws.Range("D2").Value = 0.0319362122164577
$ws.Range("F2").Value = 1.458825228356432
$ws.Range("G2").Value = 0.00242508982862988
$ws.Range("I2").Value = 0.9498458055058308
$ws.Range("M2").Value = 0.8896129614439019
$ws.Range("B3").Value = 0.5790475637801933
$ws.Range("C3").Value = 0.2242670765544972
$ws.Range("D3").Value = 0.03139188200895404
$ws.Range("F3").Value = 1.387923714114891
$ws.Range("G3").Value = 0.002430377603421512
$ws.Range("I3").Value = 0.9160288785624573
$ws.Range("M3").Value = 0.7916966004980708
$ws.Range("B4").Value = 0.5362934349374768
$ws.Range("C4").Value = 0.2046309756737514
$ws.Range("D4").Value = 0.03106734738516792
$ws.Range("F4").Value = 1.345228952570807
$ws.Range("G4").Value = 0.002433790573904417
$ws.Range("I4").Value = 0.8957553043151734
$ws.Range("M4").Value = 0.7321389162060825
$ws.Range("B5").Value = 0.5189418616138539
$ws.Range("C5").Value = 0.1966427771276926
$ws.Range("D5").Value = 0.03093745156037997
$ws.Range("F5").Value = 1.328038231075382
$ws.Range("G5").Value = 0.00243522335138507
$ws.Range("I5").Value = 0.8876149453819551
$ws.Range("M5").Value = 0.7080032487233154
$ws.Range("B6").Value = 0.5160649302519573
$ws.Range("C6").Value = 0.1953171514609835
$ws.Range("D6").Value = 0.03091602231111423
$ws.Range("F6").Value = 1.32519617743111
$ws.Range("G6").Value = 0.002435463802196808
$ws.Range("I6").Value = 0.8862705131405733
$ws.Range("M6").Value = 0.7040034629326613
$ws.Range("B7").Value = 0.5360591381586914
$ws.Range("C7").Value = 0.2045231894556991
$ws.Range("D7").Value = 0.03106558613731636
$ws.Range("F7").Value = 1.344996275626286
$ws.Range("G7").Value = 0.002433809726722164
$ws.Range("I7").Value = 0.8956450320898028
$ws.Range("M7").Value = 0.731812878540552
$ws.Range("B8").Value = 0.6249141059284398
$ws.Range("C8").Value = 0.2452690384210712
$ws.Range("D8").Value = 0.03174646773900136
$ws.Range("F8").Value = 1.434202257125122
$ws.Range("G8").Value = 0.002426878642401974
$ws.Range("I8").Value = 0.9380827999556516
$ws.Range("M8").Value = 0.8557301382571012
$ws.Range("B9").Value = 0.8014374442337839
$ws.Range("C9").Value = 0.3256529522031713
$ws.Range("D9").Value = 0.03316224863583273
$ws.Range("F9").Value = 1.615948975092238
$ws.Range("G9").Value = 0.002414598732728137
$ws.Range("I9").Value = 1.025285740100443
$ws.Range("M9").Value = 1.103538541354439
$ws.Range("B10").Value = 0.9325697272616367
$ws.Range("C10").Value = 0.3850541090577622
$ws.Range("D10").Value = 0.0342567826000959
$ws.Range("F10").Value = 1.753863709895739
$ws.Range("G10").Value = 0.002406366295235886
$ws.Range("I10").Value = 1.091917921544919
$ws.Range("M10").Value = 1.289026284251648
$ws.Range("B11").Value = 0.9925504919585819
$ws.Range("C11").Value = 0.4121630147251949
$ws.Range("D11").Value = 0.03476766346062732
$ws.Range("F11").Value = 1.817608210739763
$ws.Range("G11").Value = 0.0024027904192467
$ws.Range("I11").Value = 1.122817455939682
$ws.Range("M11").Value = 1.37426359570712
$ws.Range("B12").Value = 1.015311572959604
$ws.Range("C12").Value = 0.4224417760098618
$ws.Range("D12").Value = 0.03496307885166772
$ws.Range("F12").Value = 1.841895240696545
$ws.Range("G12").Value = 0.002401460476110574
$ws.Range("I12").Value = 1.134605222224295
$ws.Range("M12").Value = 1.40667313295836
$ws.Range("B13").Value = 1.010407438636321
$ws.Range("C13").Value = 0.4202274622537061
$ws.Range("D13").Value = 0.03492090420854765
$ws.Range("F13").Value = 1.836657936731001
$ws.Range("G13").Value = 0.002401745830978404
$ws.Range("I13").Value = 1.132062620864971
$ws.Range("M13").Value = 1.399687147544
$ws.Range("B14").Value = 0.9944221015036305
$ws.Range("C14").Value = 0.4130083858811417
$ws.Range("D14").Value = 0.03478370068155101
$ws.Range("F14").Value = 1.819603324644902
$ws.Range("G14").Value = 0.002402680520522907
$ws.Range("I14").Value = 1.123785490622112
$ws.Range("M14").Value = 1.376927251647544
$ws.Range("B15").Value = 0.9846368439373236
$ws.Range("C15").Value = 0.4085882339687146
$ws.Range("D15").Value = 0.0346999169799318
$ws.Range("F15").Value = 1.809176313793131
$ws.Range("G15").Value = 0.002403256187242939
$ws.Range("I15").Value = 1.118726876749193
$ws.Range("M15").Value = 1.363003611406384
$ws.Range("B16").Value = 0.9286564857363828
$ws.Range("C16").Value = 0.3832842883189755
$ws.Range("D16").Value = 0.03422366407912847
$ws.Range("F16").Value = 1.749718419245596
$ws.Range("G16").Value = 0.002406603376993238
$ws.Range("I16").Value = 1.089910588106719
$ws.Range("M16").Value = 1.283473823434889
$ws.Range("B17").Value = 0.8943987085215213
$ws.Range("C17").Value = 0.3677839138009062
$ws.Range("D17").Value = 0.03393488486436524
$ws.Range("F17").Value = 1.71350327274456
$ws.Range("G17").Value = 0.002408699971513433
$ws.Range("I17").Value = 1.072384941680241
$ws.Range("M17").Value = 1.234910697918394
$ws.Range("B18").Value = 0.8747253852023391
$ws.Range("C18").Value = 0.3588766567464177
$ws.Range("D18").Value = 0.03377000219592929
$ws.Range("F18").Value = 1.692767783897096
$ws.Range("G18").Value = 0.002409921802619006
$ws.Range("I18").Value = 1.062359886628755
$ws.Range("M18").Value = 1.207058700765543
$ws.Range("B19").Value = 0.868069618879872
$ws.Range("C19").Value = 0.3558621911921023
$ws.Range("D19").Value = 0.03371438177760666
$ws.Range("F19").Value = 1.68576322456039
$ws.Range("G19").Value = 0.002410338233624928
$ws.Range("I19").Value = 1.058974995041723
$ws.Range("M19").Value = 1.197642013270254
$ws.Range("B20").Value = 0.8980423116567522
$ws.Range("C20").Value = 0.3694331068298311
$ws.Range("D20").Value = 0.03396549945171756
$ws.Range("F20").Value = 1.717348629643453
$ws.Range("G20").Value = 0.002408475138017158
$ws.Range("I20").Value = 1.074244844392453
$ws.Range("M20").Value = 1.240071960443842
$ws.Range("B21").Value = 0.9991160847096126
$ws.Range("C21").Value = 0.415128439560533
$ws.Range("D21").Value = 0.0348239468736935
$ws.Range("F21").Value = 1.824608621532462
$ws.Range("G21").Value = 0.002402405325162706
$ws.Range("I21").Value = 1.126214311500505
$ws.Range("M21").Value = 1.383608731796897
$ws.Range("B22").Value = 1.065451733749114
$ws.Range("C22").Value = 0.4450702954748635
$ws.Range("D22").Value = 0.03539643513592239
$ws.Range("F22").Value = 1.895575775816525
$ws.Range("G22").Value = 0.002398579124436822
$ws.Range("I22").Value = 1.160686239952199
$ws.Range("M22").Value = 1.478191996386556
$ws.Range("B23").Value = 1.030021530707245
$ws.Range("C23").Value = 0.4290824557737665
$ws.Range("D23").Value = 0.03508981012095802
$ws.Range("F23").Value = 1.857618761235045
$ws.Range("G23").Value = 0.002400608409483497
$ws.Range("I23").Value = 1.142240810763028
$ws.Range("M23").Value = 1.427637416224755
$ws.Range("B24").Value = 0.8963949704115066
$ws.Range("C24").Value = 0.3686874938974825
$ws.Range("D24").Value = 0.03395165505232711
$ws.Range("F24").Value = 1.715609879055734
$ws.Range("G24").Value = 0.002408576733875956
$ws.Range("I24").Value = 1.073403824708549
$ws.Range("M24").Value = 1.237738343776726
$ws.Range("B25").Value = 0.7534340465837204
$ws.Range("C25").Value = 0.3038503625796523
$ws.Range("D25").Value = 0.03277007790540409
$ws.Range("F25").Value = 1.566028411643202
$ws.Range("G25").Value = 0.002417781375262405
$ws.Range("I25").Value = 1.025285740100443
$ws.Range("M25").Value = 1.035935060772204

Write-Host "Updated 168 cells for Case_3_86 380 kV case"
